# Update odds values in Sheet1 for rows 3 and 5 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 updates ---
$ws.Range("G3").Value  = 3.5
$ws.Range("I3").Value  = 2.15
$ws.Range("J3").Value  = 4.33
$ws.Range("L3").Value  = 3
$ws.Range("Q3").Value  = 1.8
$ws.Range("R3").Value  = 2.05
$ws.Range("S3").Value  = 2.35
$ws.Range("T3").Value  = 1.57
$ws.Range("AA3").Value = 2.1
$ws.Range("AB3").Value = 1.67
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 17
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 41
$ws.Range("AG3").Value = 34
$ws.Range("AI3").Value = 7
$ws.Range("AK3").Value = 19
$ws.Range("AM3").Value = 6
$ws.Range("AN3").Value = 9
$ws.Range("AO3").Value = 9.5
$ws.Range("AP3").Value = 19

# --- Row 5 updates ---
$ws.Range("G5").Value  = 2.05
$ws.Range("M5").Value  = 1.06
$ws.Range("N5").Value  = 10
$ws.Range("S5").Value  = 2.1
$ws.Range("T5").Value  = 1.7
$ws.Range("W5").Value  = 3.75
$ws.Range("X5").Value  = 1.25
$ws.Range("AA5").Value = 1.91
$ws.Range("AB5").Value = 1.8
$ws.Range("AC5").Value = 6.5
$ws.Range("AF5").Value = 19
$ws.Range("AG5").Value = 19
$ws.Range("AI5").Value = 8.5
$ws.Range("AM5").Value = 9
